# Auto-generated edit script applying the Unicorn_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 785
$ws.Range("I33").Value = 533.0952
$ws.Range("J33").Value = 1162.8572
$ws.Range("K33").Value = 533.0952
$ws.Range("L33").Value = 1162.8572
$ws.Range("M33").Value = -304.0952
$ws.Range("N33").Value = -1620.8572

$ws.Range("H34").Value = 4857.25
$ws.Range("I34").Value = 3190
$ws.Range("J34").Value = 6524.5
$ws.Range("K34").Value = 3190
$ws.Range("L34").Value = 6524.5
$ws.Range("M34").Value = -2987
$ws.Range("N34").Value = -6930.5

$ws.Range("H36").Value = 4857.25
$ws.Range("I36").Value = 3190
$ws.Range("J36").Value = 6524.5
$ws.Range("K36").Value = 3190
$ws.Range("L36").Value = 6524.5
$ws.Range("M36").Value = -2475
$ws.Range("N36").Value = -7954.5

$ws.Range("H40").Value = 1558.5
$ws.Range("I40").Value = 1402.9032
$ws.Range("J40").Value = 2247.5715
$ws.Range("K40").Value = 1402.9032
$ws.Range("L40").Value = 2247.5715
$ws.Range("M40").Value = -1227.9032
$ws.Range("N40").Value = -2597.5715

$ws.Range("H74").Value = 3615.8286
$ws.Range("I74").Value = 3177
$ws.Range("J74").Value = 4200.933
$ws.Range("K74").Value = 3177
$ws.Range("L74").Value = 4200.933
$ws.Range("M74").Value = -2241
$ws.Range("N74").Value = -6072.933

$ws.Range("H77").Value = 3615.8286
$ws.Range("I77").Value = 3177
$ws.Range("J77").Value = 4200.933
$ws.Range("K77").Value = 15885
$ws.Range("L77").Value = 21004.665
$ws.Range("M77").Value = -11205
$ws.Range("N77").Value = -30364.665

$ws.Range("H113").Value = 3748.739
$ws.Range("I113").Value = 2703.1
$ws.Range("K113").Value = 2703.1
$ws.Range("M113").Value = 550.9000000000001

$ws.Range("H132").Value = 5878.6206
$ws.Range("I132").Value = 3743.3333
$ws.Range("J132").Value = 9372.727999999999
$ws.Range("K132").Value = 11229.9999
$ws.Range("L132").Value = 28118.184
$ws.Range("M132").Value = -8699.999899999999
$ws.Range("N132").Value = -33178.18399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 550.6667
$ws.Range("I2").Value = 484.7143
$ws.Range("K2").Value = 484.7143
$ws.Range("M2").Value = -371.7143

$ws.Range("H32").Value = 7282.966
$ws.Range("I32").Value = 4539.3896
$ws.Range("K32").Value = 4539.3896
$ws.Range("M32").Value = -4252.3896

$ws.Range("H92").Value = 30550
$ws.Range("J92").Value = 30550
$ws.Range("L92").Value = 30550
$ws.Range("N92").Value = -35542

$ws.Range("H96").Value = 22500
$ws.Range("J96").Value = 22500
$ws.Range("L96").Value = 22500
$ws.Range("N96").Value = -27992

$ws.Range("H116").Value = 550.6667
$ws.Range("I116").Value = 484.7143
$ws.Range("K116").Value = 484.7143
$ws.Range("M116").Value = 1809.2857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 550.6667
$ws.Range("I3").Value = 484.7143
$ws.Range("K3").Value = 484.7143
$ws.Range("M3").Value = -370.7143

$ws.Range("H126").Value = 42500
$ws.Range("J126").Value = 42500
$ws.Range("L126").Value = 42500
$ws.Range("N126").Value = -52380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3221.5227
$ws.Range("I31").Value = 2047.6538
$ws.Range("J31").Value = 4917.1113
$ws.Range("K31").Value = 2047.6538
$ws.Range("L31").Value = 4917.1113
$ws.Range("M31").Value = -1752.6538
$ws.Range("N31").Value = -5507.1113

$ws.Range("H34").Value = 3221.5227
$ws.Range("I34").Value = 2047.6538
$ws.Range("J34").Value = 4917.1113
$ws.Range("K34").Value = 2047.6538
$ws.Range("L34").Value = 4917.1113
$ws.Range("M34").Value = -1845.6538
$ws.Range("N34").Value = -5321.1113

$ws.Range("H44").Value = 4588
$ws.Range("I44").Value = 1764
$ws.Range("J44").Value = 6000
$ws.Range("K44").Value = 1764
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = -1322
$ws.Range("N44").Value = -6884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52
$ws.Range("I12").Value = 18.5
$ws.Range("J12").Value = 85.5
$ws.Range("K12").Value = 55.5
$ws.Range("L12").Value = 256.5
$ws.Range("M12").Value = 117.5
$ws.Range("N12").Value = -602.5

$ws.Range("H131").Value = 840.95654
$ws.Range("I131").Value = 438.2857
$ws.Range("J131").Value = 913.2308
$ws.Range("K131").Value = 1314.8571
$ws.Range("L131").Value = 2739.6924
$ws.Range("M131").Value = 3725.1429
$ws.Range("N131").Value = -12819.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1282.8572
$ws.Range("I11").Value = 725
$ws.Range("J11").Value = 2026.6666
$ws.Range("K11").Value = 725
$ws.Range("L11").Value = 2026.6666
$ws.Range("M11").Value = -586
$ws.Range("N11").Value = -2304.6666

$ws.Range("H12").Value = 875.13513
$ws.Range("I12").Value = 872.9032
$ws.Range("J12").Value = 886.6667
$ws.Range("K12").Value = 872.9032
$ws.Range("L12").Value = 886.6667
$ws.Range("M12").Value = -732.9032
$ws.Range("N12").Value = -1166.6667

$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5504

$ws.Range("H36").Value = 4789.4287
$ws.Range("I36").Value = 1258.5
$ws.Range("J36").Value = 6201.8
$ws.Range("K36").Value = 1258.5
$ws.Range("L36").Value = 6201.8
$ws.Range("M36").Value = -773.5
$ws.Range("N36").Value = -7171.8

$ws.Range("H40").Value = 7000
$ws.Range("J40").Value = 7000
$ws.Range("L40").Value = 7000
$ws.Range("N40").Value = -7302

$ws.Range("H126").Value = 2573.6296
$ws.Range("I126").Value = 2411.0557
$ws.Range("J126").Value = 2898.7778
$ws.Range("K126").Value = 7233.1671
$ws.Range("L126").Value = 8696.3334
$ws.Range("M126").Value = -4763.1671
$ws.Range("N126").Value = -13636.3334

$ws.Range("H132").Value = 5055.407
$ws.Range("I132").Value = 6668.6665
$ws.Range("J132").Value = 3764.8
$ws.Range("K132").Value = 20005.9995
$ws.Range("L132").Value = 11294.4
$ws.Range("M132").Value = -17475.9995
$ws.Range("N132").Value = -16354.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 349
$ws.Range("I22").Value = 276.66666
$ws.Range("J22").Value = 457.5
$ws.Range("K22").Value = 276.66666
$ws.Range("L22").Value = 457.5
$ws.Range("M22").Value = 18.33334000000002
$ws.Range("N22").Value = -1047.5

$ws.Range("H27").Value = 349
$ws.Range("I27").Value = 276.66666
$ws.Range("J27").Value = 457.5
$ws.Range("K27").Value = 276.66666
$ws.Range("L27").Value = 457.5
$ws.Range("M27").Value = -169.66666
$ws.Range("N27").Value = -671.5

$ws.Range("H33").Value = 6700
$ws.Range("I33").Value = 3750
$ws.Range("J33").Value = 8666.666999999999
$ws.Range("K33").Value = 3750
$ws.Range("L33").Value = 8666.666999999999
$ws.Range("M33").Value = -3460
$ws.Range("N33").Value = -9246.666999999999

$ws.Range("H136").Value = 5227.9736
$ws.Range("I136").Value = 3123.7
$ws.Range("J136").Value = 7566.0557
$ws.Range("K136").Value = 9371.099999999999
$ws.Range("L136").Value = 22698.1671
$ws.Range("M136").Value = -6821.099999999999
$ws.Range("N136").Value = -27798.1671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 155.5238
$ws.Range("I113").Value = 155.5238
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 466.5714
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1703.4286
$ws.Range("N113").ClearContents()
